$wb = $excel.ActiveWorkbook

# This script applies updated market-price figures (columns H-N) that were
# refreshed by the scheduled market-data runner, across the ALC, ARM, CRP,
# CUL, GSM, LTW and WVR sheets. Values are written directly into the cells
# backing each sheet's Table_XXX listobject.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 97.72727
$ws.Range("I9").Value = 57.8
$ws.Range("K9").Value = 57.8
$ws.Range("M9").Value = 111.2
$ws.Range("H17").Value = 1617.3928
$ws.Range("J17").Value = 1794.8636
$ws.Range("L17").Value = 5384.5908
$ws.Range("N17").Value = -5720.5908
$ws.Range("H80").Value = 704.6429000000001
$ws.Range("I80").Value = 709.375
$ws.Range("K80").Value = 2128.125
$ws.Range("M80").Value = -1130.125
$ws.Range("H83").Value = 704.6429000000001
$ws.Range("I83").Value = 709.375
$ws.Range("K83").Value = 6384.375
$ws.Range("M83").Value = -1392.375
$ws.Range("H86").Value = 4500
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("H89").Value = 4500
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384
$ws.Range("H131").Value = 1509.8182
$ws.Range("J131").Value = 4490
$ws.Range("L131").Value = 13470
$ws.Range("N131").Value = -23550
$ws.Range("H132").Value = 14919.167
$ws.Range("I132").Value = 14410
$ws.Range("K132").Value = 43230
$ws.Range("M132").Value = -40700
$ws.Range("H135").Value = 695.3
$ws.Range("I135").Value = 630.3182
$ws.Range("K135").Value = 5672.8638
$ws.Range("M135").Value = -3137.8638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1587.5714
$ws.Range("I26").Value = 942.6
$ws.Range("J26").Value = 3200
$ws.Range("K26").Value = 942.6
$ws.Range("L26").Value = 3200
$ws.Range("M26").Value = -612.6
$ws.Range("N26").Value = -3860
$ws.Range("H61").Value = 3127.375
$ws.Range("I61").Value = 3023.8
$ws.Range("K61").Value = 3023.8
$ws.Range("M61").Value = -2811.8
$ws.Range("H136").Value = 3127.375
$ws.Range("I136").Value = 3023.8
$ws.Range("K136").Value = 9071.400000000001
$ws.Range("M136").Value = -6521.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 874.25
$ws.Range("J22").Value = 848.75
$ws.Range("L22").Value = 848.75
$ws.Range("N22").Value = -1548.75
$ws.Range("H33").Value = 5133.222
$ws.Range("I33").Value = 1866.5
$ws.Range("K33").Value = 1866.5
$ws.Range("M33").Value = -1487.5
$ws.Range("H44").Value = 14999.333
$ws.Range("I44").Value = 5000
$ws.Range("K44").Value = 5000
$ws.Range("M44").Value = -4558
$ws.Range("H55").Value = 4649.75
$ws.Range("I55").Value = 4649.75
$ws.Range("K55").Value = 4649.75
$ws.Range("M55").Value = -4334.75
$ws.Range("H58").Value = 2949.348
$ws.Range("I58").Value = 1729.2106
$ws.Range("K58").Value = 1729.2106
$ws.Range("M58").Value = -1526.2106
$ws.Range("H68").Value = 39688.75
$ws.Range("J68").Value = 59377.5
$ws.Range("L68").Value = 59377.5
$ws.Range("N68").Value = -60875.5
$ws.Range("H71").Value = 39688.75
$ws.Range("J71").Value = 59377.5
$ws.Range("L71").Value = 178132.5
$ws.Range("N71").Value = -185620.5
$ws.Range("H86").Value = 6083.3335
$ws.Range("I86").Value = 6083.3335
$ws.Range("K86").Value = 6083.3335
$ws.Range("M86").Value = -4960.3335
$ws.Range("H89").Value = 6083.3335
$ws.Range("I89").Value = 6083.3335
$ws.Range("K89").Value = 30416.6675
$ws.Range("M89").Value = -24800.6675
$ws.Range("H136").Value = 2949.348
$ws.Range("I136").Value = 1729.2106
$ws.Range("K136").Value = 5187.6318
$ws.Range("M136").Value = -2637.6318

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 722.3
$ws.Range("I44").Value = 203.83333
$ws.Range("J44").Value = 1500
$ws.Range("K44").Value = 611.49999
$ws.Range("L44").Value = 4500
$ws.Range("M44").Value = -213.49999
$ws.Range("N44").Value = -5296
$ws.Range("H107").Value = 428.66666
$ws.Range("I107").Value = 492
$ws.Range("J107").Value = 397
$ws.Range("K107").Value = 1476
$ws.Range("L107").Value = 1191
$ws.Range("M107").Value = 444
$ws.Range("N107").Value = -5031
$ws.Range("H127").Value = 2000
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920
$ws.Range("H131").Value = 2145.2942
$ws.Range("J131").Value = 2731.889
$ws.Range("L131").Value = 8195.667000000001
$ws.Range("N131").Value = -18275.667
$ws.Range("H132").Value = 3666.6667
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 33000.0003
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -38060.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4534.6665
$ws.Range("I113").Value = 6202
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 6202
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = -4032
$ws.Range("N113").Value = -5540
$ws.Range("H122").Value = 2426.5
$ws.Range("I122").Value = 1568.6666
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4705.9998
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2255.9998
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 31167.277
$ws.Range("I132").Value = 35484.773
$ws.Range("K132").Value = 106454.319
$ws.Range("M132").Value = -103924.319

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2043.4
$ws.Range("J46").Value = 998.3333
$ws.Range("L46").Value = 998.3333
$ws.Range("N46").Value = -1374.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1331.1428
$ws.Range("J113").Value = 1321.7778
$ws.Range("L113").Value = 3965.3334
$ws.Range("N113").Value = -8305.3334
$ws.Range("H122").Value = 3026.6667
$ws.Range("I122").Value = 2335.5454
$ws.Range("J122").Value = 4927.25
$ws.Range("K122").Value = 7006.6362
$ws.Range("L122").Value = 14781.75
$ws.Range("M122").Value = -4556.6362
$ws.Range("N122").Value = -19681.75
$ws.Range("H132").Value = 2383.6155
$ws.Range("I132").Value = 1899.8889
$ws.Range("J132").Value = 3472
$ws.Range("K132").Value = 5699.6667
$ws.Range("L132").Value = 10416
$ws.Range("M132").Value = -3169.6667
$ws.Range("N132").Value = -15476
